# ============================================================================
# Applies the "gh-pages output generated at 456a3b4" update to
# 北京-漫展信息.xlsx:
#   - Sheet "展览"   (index 1): numeric refreshes + 2 new rows (37 & 41)
#   - Sheet "演出"   (index 2): numeric refreshes only
#   - Sheet "本地生活"(index 3): untouched
#   - Sheet "全部类型"(index 4): numeric refreshes + 2 new rows (43 & 48)
#     (mirrors sheet "展览", shifted down by 4 rows because it also holds
#     the "演出" + "本地生活" rows ahead of the "展览" rows)
# ============================================================================

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)   # 展览
$ws2 = $wb.Worksheets.Item(2)   # 演出
$ws4 = $wb.Worksheets.Item(4)   # 全部类型

# ----------------------------------------------------------------------------
# Helper: insert one brand-new data row into $ws at row $r, re-using the
# formatting of the row right above it ($r-1) so the inserted row ends up
# with the same look (incl. the bold/bordered index-column style) as every
# other data row, and force the text columns to Text format first so Excel's
# autodetection doesn't turn the "YYYY-MM-DD" strings / urls into dates or
# numbers.
# ----------------------------------------------------------------------------
function Insert-Row($ws, $r, $a, $b, $c, $d, $e, $f, $g, $h, $i) {
    $ws.Rows.Item($r).Insert()

    $srcRow = $r - 1

    foreach ($col in @("A","B","C","D","E","F","G","H","I")) {
        $ws.Range("$col$srcRow").Copy()
        $ws.Range("$col$r").PasteSpecial(-4122)  # xlPasteFormats
    }

    foreach ($col in @("B","C","D","E","H","I")) {
        $ws.Range("$col$r").NumberFormat = "@"
    }

    $ws.Range("A$r").Value = $a
    $ws.Range("B$r").Value = $b
    $ws.Range("C$r").Value = $c
    $ws.Range("D$r").Value = $d
    $ws.Range("E$r").Value = $e
    $ws.Range("F$r").Value = $f
    $ws.Range("G$r").Value = $g
    $ws.Range("H$r").Value = $h
    $ws.Range("I$r").Value = $i

    # Re-apply plain-text formatting from the (now text-valued) neighbour row
    # so the inserted cells don't keep a stray explicit "@" number format
    # that the original workbook never used.
    foreach ($col in @("B","C","D","E","H","I")) {
        $ws.Range("$col$srcRow").Copy()
        $ws.Range("$col$r").PasteSpecial(-4122)  # xlPasteFormats
    }
}

# ----------------------------------------------------------------------------
# Sheet "展览" (ws1) - numeric refreshes
# ----------------------------------------------------------------------------
$ws1.Range("F4").Value = 3398
$ws1.Range("I4").Value = "//i0.hdslb.com/bfs/openplatform/202402/MQV2WfTg1707199914411.jpeg"
$ws1.Range("F5").Value = 224
$ws1.Range("F6").Value = 4895
$ws1.Range("F7").Value = 481
$ws1.Range("F8").Value = 311
$ws1.Range("F10").Value = 643
$ws1.Range("F11").Value = 288
$ws1.Range("F12").Value = 52
$ws1.Range("F13").Value = 19
$ws1.Range("F14").Value = 673
$ws1.Range("F15").Value = 291
$ws1.Range("F16").Value = 27
$ws1.Range("F20").Value = 4788
$ws1.Range("F21").Value = 30
$ws1.Range("F24").Value = 5921
$ws1.Range("F26").Value = 2
$ws1.Range("F27").Value = 1805
$ws1.Range("F28").Value = 252
$ws1.Range("F29").Value = 682
$ws1.Range("F30").Value = 4429
$ws1.Range("F31").Value = 308
$ws1.Range("F32").Value = 101
$ws1.Range("F33").Value = 133
$ws1.Range("F34").Value = 888
$ws1.Range("F36").Value = 15

# Row 37 ("北京·IDO动漫游戏嘉年华45th") gets pushed to row 38 with an updated
# 想去人数 (806 -> 812); row 38 ("北京·第16届IJOY...") becomes row 39 with its
# own refresh (878 -> 884); row 39 ("北京·原神only") becomes row 40 with its
# refresh (6 -> 7). A brand new row is inserted at 37, and another one is
# appended at the (new) end, row 41.
Insert-Row $ws1 37 36 "2024-04-19" "北京·次元风暴游园会" "安翔路1号院 老故事503文化创意产业园" "2024.04.19 09:00-04.20 17:00" 0 55 "https://show.bilibili.com/platform/detail.html?id=81781" "//i0.hdslb.com/bfs/openplatform/202402/QsBPojEU1707191707677.jpeg"

$ws1.Range("F38").Value = 812
$ws1.Range("F39").Value = 884
$ws1.Range("F40").Value = 7

# Row 41 is a brand new row appended right after the current last row (40),
# so a plain Insert (which would push non-existent content) isn't needed -
# Insert-Row still works here since Rows.Item(41).Insert() on the first
# unused row is a no-op shift and leaves us a clean row to fill in.
Insert-Row $ws1 41 40 "2024-05-18" "北京·次元风暴游园会2.0" "安翔路1号院 老故事503文化创意产业园" "2024.05.18 09:00-05.19 17:00" 0 55 "https://show.bilibili.com/platform/detail.html?id=81782" "//i0.hdslb.com/bfs/openplatform/202402/7J276vFp1707191576670.jpeg"

# ----------------------------------------------------------------------------
# Sheet "演出" (ws2) - numeric refreshes only
# ----------------------------------------------------------------------------
$ws2.Range("F2").Value = 7
$ws2.Range("F6").Value = 49

# ----------------------------------------------------------------------------
# Sheet "全部类型" (ws4) - numeric refreshes
# ----------------------------------------------------------------------------
$ws4.Range("F7").Value = 7
$ws4.Range("F8").Value = 3398
$ws4.Range("I8").Value = "//i0.hdslb.com/bfs/openplatform/202402/MQV2WfTg1707199914411.jpeg"
$ws4.Range("F9").Value = 224
$ws4.Range("F10").Value = 4895
$ws4.Range("F11").Value = 481
$ws4.Range("F12").Value = 312
$ws4.Range("F14").Value = 643
$ws4.Range("F15").Value = 288
$ws4.Range("F16").Value = 52
$ws4.Range("F17").Value = 19
$ws4.Range("F18").Value = 673
$ws4.Range("F19").Value = 291
$ws4.Range("F20").Value = 27
$ws4.Range("F25").Value = 4788
$ws4.Range("F26").Value = 30
$ws4.Range("F29").Value = 5921
$ws4.Range("F31").Value = 2
$ws4.Range("F32").Value = 1808
$ws4.Range("F33").Value = 252
$ws4.Range("F34").Value = 682
$ws4.Range("F35").Value = 4429
$ws4.Range("F36").Value = 308
$ws4.Range("F38").Value = 101
$ws4.Range("F39").Value = 133
$ws4.Range("F40").Value = 888
$ws4.Range("F42").Value = 15

Insert-Row $ws4 43 42 "2024-04-19" "北京·次元风暴游园会" "安翔路1号院 老故事503文化创意产业园" "2024.04.19 09:00-04.20 17:00" 0 55 "https://show.bilibili.com/platform/detail.html?id=81781" "//i0.hdslb.com/bfs/openplatform/202402/QsBPojEU1707191707677.jpeg"

$ws4.Range("F44").Value = 812
$ws4.Range("F45").Value = 884
$ws4.Range("F47").Value = 7

Insert-Row $ws4 48 47 "2024-05-18" "北京·次元风暴游园会2.0" "安翔路1号院 老故事503文化创意产业园" "2024.05.18 09:00-05.19 17:00" 0 55 "https://show.bilibili.com/platform/detail.html?id=81782" "//i0.hdslb.com/bfs/openplatform/202402/7J276vFp1707191576670.jpeg"

Write-Host "edit complete"
